# Generate Report for handoff
# - Renames the UUID-named markdown/xlf artifacts (new localization run)
# - Updates the handoff timestamps for zh-cn / de-de
# - Drops the row for the file that previously failed transform
#   (f5864e5f-6d52-4c13-892b-d6126c317490.md / "Handoff transform failed"),
#   which no longer appears in the report; the ".localization-config" row
#   shifts up to take its place.

$wb = $excel.ActiveWorkbook

$newGuid = "2c026530-0b0e-4853-a10c-c1820915f1f6"
$newHash = "9ffe65e21385c251ef3e51ba71d480539f8c25da"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Remove the "Handoff transform failed" row (row 3); ".localization-config"
# row shifts up from row 4 to row 3.
$wsOverview.Rows.Item(3).Delete()

# Rename the artifact referenced in row 2.
$wsOverview.Range("A2").Value = $newGuid + ".md"

# Rebuild hyperlinks so relationship ids stay sequential / match the new
# layout (old hyperlink collection still references the removed row).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e2725f2f080e7313645aaa1efd780ff4d548714/e2e/" + $newGuid + ".md", [Type]::Missing, [Type]::Missing, $newGuid + ".md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e2725f2f080e7313645aaa1efd780ff4d548714/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(3).Delete()

$wsZh.Range("A2").Value = $newGuid + ".md"
$wsZh.Range("C2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-02-15 08:39:19"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e2725f2f080e7313645aaa1efd780ff4d548714/e2e/" + $newGuid + ".md", [Type]::Missing, [Type]::Missing, $newGuid + ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e30301d7161720e69171a43a3025ac1c5f21aeb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $newGuid + "." + $newHash + ".zh-cn.xlf", [Type]::Missing, [Type]::Missing, $newGuid + "." + $newHash + ".zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e2725f2f080e7313645aaa1efd780ff4d548714/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(3).Delete()

$wsDe.Range("A2").Value = $newGuid + ".md"
$wsDe.Range("C2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDe.Range("D2").Value = "2016-02-15 08:39:33"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e2725f2f080e7313645aaa1efd780ff4d548714/e2e/" + $newGuid + ".md", [Type]::Missing, [Type]::Missing, $newGuid + ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a2e1070a3c6d3d39c923eb2babe4fff84762031/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $newGuid + "." + $newHash + ".de-de.xlf", [Type]::Missing, [Type]::Missing, $newGuid + "." + $newHash + ".de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e2725f2f080e7313645aaa1efd780ff4d548714/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
